# Restored from revision #a062a2943b35dd155eccb856f56c65c54b2caa07.TEST Author: admin. Type: SAVE.
# Update the "Integer min" threshold for rule R30 (row 10) on the Rules
# sheet: cell C10 changes from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")
$ws.Range("C10").Value = 1
